$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new values look numeric (e.g. "1.006"),
# so Excel keeps them as text strings instead of converting to numbers.
$textCells = @("D5","D7","D8","D9","D10","D11","D13","D15","D17","D18","D19","D20","D22","D24","D26","D27","D28","D29","D30","D31","D32","D33","D35","D36","D37","D38","D39","D40","D41","D43","D44","D45","D46","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values scraped for this run
$ws.Range("D2").Value2 = '28.639.65'
$ws.Range("E2").Value2 = '  +1.19%  '
$ws.Range("D3").Value2 = '1.867.50'
$ws.Range("E3").Value2 = '  +1.43%  '
$ws.Range("E4").Value2 = '  +0.13%  '
$ws.Range("D5").Value2 = '326.58'
$ws.Range("E5").Value2 = '  -2.11%  '
$ws.Range("D7").Value2 = '0.4656'
$ws.Range("E7").Value2 = '  +0.79%  '
$ws.Range("D8").Value2 = '0.3917'
$ws.Range("E8").Value2 = '  +0.99%  '
$ws.Range("D9").Value2 = '0.07903'
$ws.Range("E9").Value2 = '  +0.31%  '
$ws.Range("D10").Value2 = '0.9743'
$ws.Range("E10").Value2 = '  +0.41%  '
$ws.Range("D11").Value2 = '22.29'
$ws.Range("E11").Value2 = '  +1.26%  '
$ws.Range("D12").Value2 = '1.850.68'
$ws.Range("E12").Value2 = '  -2.03%  '
$ws.Range("D13").Value2 = '5.740'
$ws.Range("E13").Value2 = '  -1.00%  '
$ws.Range("E14").Value2 = '  +0.14%  '
$ws.Range("D15").Value2 = '0.06915'
$ws.Range("E15").Value2 = '  +0.04%  '
$ws.Range("E16").Value2 = '  +1.63%  '
$ws.Range("D17").Value2 = '1.006'
$ws.Range("E17").Value2 = '  +0.08%  '
$ws.Range("D18").Value2 = '0.00001004'
$ws.Range("E18").Value2 = '  +0.63%  '
$ws.Range("D19").Value2 = '16.94'
$ws.Range("E19").Value2 = '  -0.26%  '
$ws.Range("D20").Value2 = '1.005'
$ws.Range("E20").Value2 = '  -0.06%  '
$ws.Range("D21").Value2 = '28.602.59'
$ws.Range("E21").Value2 = '  +0.95%  '
$ws.Range("D22").Value2 = '5.332'
$ws.Range("E22").Value2 = '  -0.57%  '
$ws.Range("E23").Value2 = '  -0.69%  '
$ws.Range("D24").Value2 = '2.126'
$ws.Range("E24").Value2 = '  -2.18%  '
$ws.Range("D25").Value2 = '2.063.27'
$ws.Range("E25").Value2 = '  -1.69%  '
$ws.Range("D26").Value2 = '155.35'
$ws.Range("E26").Value2 = '  +1.15%  '
$ws.Range("D27").Value2 = '19.30'
$ws.Range("E27").Value2 = '  -0.18%  '
$ws.Range("D28").Value2 = '5.780'
$ws.Range("E28").Value2 = '  -2.55%  '
$ws.Range("D29").Value2 = '1.993'
$ws.Range("E29").Value2 = '  +0.46%  '
$ws.Range("D30").Value2 = '119.31'
$ws.Range("E30").Value2 = '  +1.84%  '
$ws.Range("D31").Value2 = '0.09365'
$ws.Range("E31").Value2 = '  +0.18%  '
$ws.Range("D32").Value2 = '0.9397'
$ws.Range("E32").Value2 = '  -1.40%  '
$ws.Range("D33").Value2 = '5.325'
$ws.Range("E33").Value2 = '  -0.14%  '
$ws.Range("E34").Value2 = '  +0.85%  '
$ws.Range("D35").Value2 = '3.348'
$ws.Range("E35").Value2 = '  -3.41%  '
$ws.Range("D36").Value2 = '0.05844'
$ws.Range("E36").Value2 = '  -3.47%  '
$ws.Range("D37").Value2 = '0.02118'
$ws.Range("E37").Value2 = '  -2.88%  '
$ws.Range("D38").Value2 = '1.157'
$ws.Range("E38").Value2 = '  +0.24%  '
$ws.Range("D39").Value2 = '7.863'
$ws.Range("E39").Value2 = '  +3.17%  '
$ws.Range("D40").Value2 = '0.5656'
$ws.Range("E40").Value2 = '  +0.00%  '
$ws.Range("D41").Value2 = '9.977'
$ws.Range("E41").Value2 = '  -0.50%  '
$ws.Range("E42").Value2 = '  -0.73%  '
$ws.Range("D43").Value2 = '0.07355'
$ws.Range("E43").Value2 = '  +4.49%  '
$ws.Range("B44").Value2 = 'RenderToken'
$ws.Range("C44").Value2 = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").Value2 = '2.199'
$ws.Range("E44").Value2 = '  -8.14%  '
$ws.Range("B45").Value2 = 'Decentraland'
$ws.Range("C45").Value2 = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").Value2 = '0.5327'
$ws.Range("E45").Value2 = '  -0.04%  '
$ws.Range("B46").Value2 = 'EnergySwap'
$ws.Range("C46").Value2 = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value2 = '11.65'
$ws.Range("E46").Value2 = '  -0.29%  '
$ws.Range("E47").Value2 = '  -7.00%  '
$ws.Range("E48").Value2 = '  -0.53%  '
$ws.Range("D49").Value2 = '114.10'
$ws.Range("E49").Value2 = '  +1.03%  '
$ws.Range("D50").Value2 = '2.356'
$ws.Range("E50").Value2 = '  +1.28%  '
$ws.Range("D51").Value2 = '1.006'
$ws.Range("E51").Value2 = '  +0.01%  '
